$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update results for Steel
$ws.Range("B3").Value = 13243.3408437119
$ws.Range("D6").Value = 53232.37434478001
